$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the localhost URL used in the Steps column of TC_001 (row 2, column D) ---
# "1. Navigate to http://localhost:3001/" -> "1. Navigate to http://localhost:3000/"
$ws.Range("D2").Value = "1. Navigate to http://localhost:3000/"

# --- Apply explicit ("best fit") column widths to B:E, matching the widths Excel
#     would compute after an AutoFit pass over the Test Case Title / Pre-conditions /
#     Steps / Expected Result columns ---
$ws.Columns.Item(2).ColumnWidth = 31.0
$ws.Columns.Item(3).ColumnWidth = 18.166666666666668
$ws.Columns.Item(4).ColumnWidth = 56.0
$ws.Columns.Item(5).ColumnWidth = 119.0

# --- Update the view: scroll the window so column D is the left-most visible
#     column, and select the entire column F (the "Automated" column) ---
$excel.ActiveWindow.ScrollColumn = 4
$ws.Columns.Item(6).Select()
